$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original cell styles, force text format while assigning values
# so Excel does not auto-convert numeric-looking strings (e.g. "210.70")
# into actual numbers, then restore the original style/format.
$priceVolRange = $ws.Range("D2:E51")
$origStyle = $priceVolRange.Style
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.117.88'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").Value = '1.667.32'
$ws.Range("E3").Value = '  -1.34%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.62%  '
$ws.Range("D5").Value = '210.70'
$ws.Range("E5").Value = '  -3.31%  '
$ws.Range("D6").Value = '0.5248'
$ws.Range("E6").Value = '  -2.32%  '
$ws.Range("D7").Value = '1.002'
$ws.Range("D8").Value = '0.2627'
$ws.Range("E8").Value = '  -3.71%  '
$ws.Range("D9").Value = '0.06287'
$ws.Range("E9").Value = '  -2.30%  '
$ws.Range("D10").Value = '21.14'
$ws.Range("E10").Value = '  -1.91%  '
$ws.Range("D11").Value = '0.07535'
$ws.Range("E11").Value = '  -1.70%  '
$ws.Range("D12").Value = '1.666.88'
$ws.Range("E12").Value = '  -1.48%  '
$ws.Range("E13").Value = '  -1.99%  '
$ws.Range("D14").Value = '0.5546'
$ws.Range("E14").Value = '  -4.02%  '
$ws.Range("D15").Value = '66.70'
$ws.Range("E15").Value = '  -0.20%  '
$ws.Range("D16").Value = '0.000007934'
$ws.Range("E16").Value = '  -5.15%  '
$ws.Range("D17").Value = '26.134.36'
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("D18").Value = '1.002'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").Value = '4.729'
$ws.Range("E19").Value = '  -3.46%  '
$ws.Range("D20").Value = '186.35'
$ws.Range("E20").Value = '  -2.15%  '
$ws.Range("D21").Value = '10.33'
$ws.Range("E21").Value = '  -4.86%  '
$ws.Range("D22").Value = '6.160'
$ws.Range("E22").Value = '  -1.48%  '
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  -0.65%  '
$ws.Range("D24").Value = '149.55'
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("D26").Value = '7.476'
$ws.Range("E26").Value = '  -4.69%  '
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("D28").Value = '0.06257'
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("D29").Value = '1.354'
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("D30").Value = '1.279'
$ws.Range("E30").Value = '  -3.50%  '
$ws.Range("D31").Value = '3.506'
$ws.Range("E31").Value = '  -2.42%  '
$ws.Range("D32").Value = '3.410'
$ws.Range("E32").Value = '  -4.74%  '
$ws.Range("D33").Value = '1.631'
$ws.Range("E33").Value = '  -2.33%  '
$ws.Range("E34").Value = '  -3.04%  '
$ws.Range("D35").Value = '2.412'
$ws.Range("E35").Value = '  -0.29%  '
$ws.Range("D36").Value = '0.6022'
$ws.Range("E36").Value = '  -2.05%  '
$ws.Range("D37").Value = '2.729'
$ws.Range("E37").Value = '  -1.17%  '
$ws.Range("D38").Value = '6.112'
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("D39").Value = '1.105.31'
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("E40").Value = '  -2.30%  '
$ws.Range("D41").Value = '0.8698'
$ws.Range("E41").Value = '  -1.18%  '
$ws.Range("E42").Value = '  -1.04%  '
$ws.Range("E43").Value = '  -1.27%  '
$ws.Range("D44").Value = '1.820.41'
$ws.Range("E44").Value = '  -1.15%  '
$ws.Range("E45").Value = '  +0.86%  '
$ws.Range("E46").Value = '  -3.86%  '
$ws.Range("D48").Value = '8.050'
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("D49").Value = '0.05232'
$ws.Range("E49").Value = '  -0.98%  '
$ws.Range("D50").Value = '0.4242'
$ws.Range("E50").Value = '  -1.34%  '
$ws.Range("D51").Value = '5.966'
$ws.Range("E51").Value = '  -1.02%  '
$priceVolRange.Style = $origStyle
